$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update string columns B, C, D (renaming BRAVER -> SOMERSAULT, HIGHWIND -> HELLMASKER) ---
for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 2).Value2 = "Missile_SOMERSAULT_State_Update"
    $ws.Cells.Item($r, 3).Value2 = "MISSILE_SOMERSAULT_413.MISSILE_SOMERSAULT_413"
    $ws.Cells.Item($r, 4).Value2 = "MISSILE_SOMERSAULT"
}
for ($r = 55; $r -le 107; $r++) {
    $ws.Cells.Item($r, 2).Value2 = "Missile_HELLMASKER_State_Update"
    $ws.Cells.Item($r, 3).Value2 = "MISSILE_HELLMASKER_309.MISSILE_HELLMASKER_309"
    $ws.Cells.Item($r, 4).Value2 = "MISSILE_HELLMASKER"
}

# --- Update numeric columns F, G, H, I, J, K ---
$ws.Cells.Item(2, 6).Value2 = 1116576.797998612
$ws.Cells.Item(2, 7).Value2 = 4841127.592487725
$ws.Cells.Item(2, 8).Value2 = 3985228.017267243
$ws.Cells.Item(2, 9).Value2 = 1114862.618608733
$ws.Cells.Item(2, 10).Value2 = 4843220.03141693
$ws.Cells.Item(2, 11).Value2 = 3984374.979855397
$ws.Cells.Item(3, 6).Value2 = 1116576.797998612
$ws.Cells.Item(3, 7).Value2 = 4841127.592487725
$ws.Cells.Item(3, 8).Value2 = 3985228.017267243
$ws.Cells.Item(3, 9).Value2 = 1114892.232638698
$ws.Cells.Item(3, 10).Value2 = 4843171.379861275
$ws.Cells.Item(3, 11).Value2 = 3984678.776202583
$ws.Cells.Item(4, 6).Value2 = 1116576.797998612
$ws.Cells.Item(4, 7).Value2 = 4841127.592487725
$ws.Cells.Item(4, 8).Value2 = 3985228.017267243
$ws.Cells.Item(4, 9).Value2 = 1114922.575886906
$ws.Cells.Item(4, 10).Value2 = 4843122.728305617
$ws.Cells.Item(4, 11).Value2 = 3984967.423871347
$ws.Cells.Item(5, 6).Value2 = 1116576.797998612
$ws.Cells.Item(5, 7).Value2 = 4841127.592487725
$ws.Cells.Item(5, 8).Value2 = 3985228.017267243
$ws.Cells.Item(5, 9).Value2 = 1114953.666309685
$ws.Cells.Item(5, 10).Value2 = 4843074.076749961
$ws.Cells.Item(5, 11).Value2 = 3985240.922861687
$ws.Cells.Item(6, 6).Value2 = 1116576.797998612
$ws.Cells.Item(6, 7).Value2 = 4841127.592487725
$ws.Cells.Item(6, 8).Value2 = 3985228.017267243
$ws.Cells.Item(6, 9).Value2 = 1114985.522305524
$ws.Cells.Item(6, 10).Value2 = 4843025.425194304
$ws.Cells.Item(6, 11).Value2 = 3985499.273173604
$ws.Cells.Item(7, 6).Value2 = 1116576.797998612
$ws.Cells.Item(7, 7).Value2 = 4841127.592487725
$ws.Cells.Item(7, 8).Value2 = 3985228.017267243
$ws.Cells.Item(7, 9).Value2 = 1115018.162725952
$ws.Cells.Item(7, 10).Value2 = 4842976.773638649
$ws.Cells.Item(7, 11).Value2 = 3985742.474807098
$ws.Cells.Item(8, 6).Value2 = 1116576.797998612
$ws.Cells.Item(8, 7).Value2 = 4841127.592487725
$ws.Cells.Item(8, 8).Value2 = 3985228.017267243
$ws.Cells.Item(8, 9).Value2 = 1115051.606886706
$ws.Cells.Item(8, 10).Value2 = 4842928.122082992
$ws.Cells.Item(8, 11).Value2 = 3985970.527762169
$ws.Cells.Item(9, 6).Value2 = 1116576.797998612
$ws.Cells.Item(9, 7).Value2 = 4841127.592487725
$ws.Cells.Item(9, 8).Value2 = 3985228.017267243
$ws.Cells.Item(9, 9).Value2 = 1115085.874579149
$ws.Cells.Item(9, 10).Value2 = 4842879.470527335
$ws.Cells.Item(9, 11).Value2 = 3986183.432038816
$ws.Cells.Item(10, 6).Value2 = 1116576.797998612
$ws.Cells.Item(10, 7).Value2 = 4841144.105317269
$ws.Cells.Item(10, 8).Value2 = 3985228.017267243
$ws.Cells.Item(10, 9).Value2 = 1115120.986081993
$ws.Cells.Item(10, 10).Value2 = 4842830.818971679
$ws.Cells.Item(10, 11).Value2 = 3986381.18763704
$ws.Cells.Item(11, 6).Value2 = 1116528.117509214
$ws.Cells.Item(11, 7).Value2 = 4841160.618146812
$ws.Cells.Item(11, 8).Value2 = 3985427.271658464
$ws.Cells.Item(11, 9).Value2 = 1115156.962173293
$ws.Cells.Item(11, 10).Value2 = 4842782.167416023
$ws.Cells.Item(11, 11).Value2 = 3986563.794556841
$ws.Cells.Item(12, 6).Value2 = 1116499.520163769
$ws.Cells.Item(12, 7).Value2 = 4841177.130976355
$ws.Cells.Item(12, 8).Value2 = 3985545.957764132
$ws.Cells.Item(12, 9).Value2 = 1115193.824142744
$ws.Cells.Item(12, 10).Value2 = 4842733.515860366
$ws.Cells.Item(12, 11).Value2 = 3986731.25279822
$ws.Cells.Item(13, 6).Value2 = 1116481.467421317
$ws.Cells.Item(13, 7).Value2 = 4841193.643805898
$ws.Cells.Item(13, 8).Value2 = 3985630.774702165
$ws.Cells.Item(13, 9).Value2 = 1115231.59380428
$ws.Cells.Item(13, 10).Value2 = 4842684.864304709
$ws.Cells.Item(13, 11).Value2 = 3986883.562361174
$ws.Cells.Item(14, 6).Value2 = 1116468.877512012
$ws.Cells.Item(14, 7).Value2 = 4841210.15663544
$ws.Cells.Item(14, 8).Value2 = 3985696.821606371
$ws.Cells.Item(14, 9).Value2 = 1115270.293508984
$ws.Cells.Item(14, 10).Value2 = 4842636.212749053
$ws.Cells.Item(14, 11).Value2 = 3987020.723245705
$ws.Cells.Item(15, 6).Value2 = 1116459.485062645
$ws.Cells.Item(15, 7).Value2 = 4841226.669464983
$ws.Cells.Item(15, 8).Value2 = 3985750.91914398
$ws.Cells.Item(15, 9).Value2 = 1115309.946158313
$ws.Cells.Item(15, 10).Value2 = 4842587.561193397
$ws.Cells.Item(15, 11).Value2 = 3987142.735451814
$ws.Cells.Item(16, 6).Value2 = 1116452.141478857
$ws.Cells.Item(16, 7).Value2 = 4841243.182294525
$ws.Cells.Item(16, 8).Value2 = 3985796.735809407
$ws.Cells.Item(16, 9).Value2 = 1115350.575217652
$ws.Cells.Item(16, 10).Value2 = 4842538.90963774
$ws.Cells.Item(16, 11).Value2 = 3987249.598979499
$ws.Cells.Item(17, 6).Value2 = 1116446.199984869
$ws.Cells.Item(17, 7).Value2 = 4841259.695124068
$ws.Cells.Item(17, 8).Value2 = 3985836.473422206
$ws.Cells.Item(17, 9).Value2 = 1115392.2047302
$ws.Cells.Item(17, 10).Value2 = 4842490.258082083
$ws.Cells.Item(17, 11).Value2 = 3987341.313828761
$ws.Cells.Item(18, 6).Value2 = 1116441.266427785
$ws.Cells.Item(18, 7).Value2 = 4841276.207953612
$ws.Cells.Item(18, 8).Value2 = 3985871.557802123
$ws.Cells.Item(18, 9).Value2 = 1115434.859331195
$ws.Cells.Item(18, 10).Value2 = 4842441.606526427
$ws.Cells.Item(18, 11).Value2 = 3987417.8799996
$ws.Cells.Item(19, 6).Value2 = 1116437.085609361
$ws.Cells.Item(19, 7).Value2 = 4841292.720783154
$ws.Cells.Item(19, 8).Value2 = 3985902.965382141
$ws.Cells.Item(19, 9).Value2 = 1115478.564262498
$ws.Cells.Item(19, 10).Value2 = 4842392.954970771
$ws.Cells.Item(19, 11).Value2 = 3987479.297492016
$ws.Cells.Item(20, 6).Value2 = 1116433.484242171
$ws.Cells.Item(20, 7).Value2 = 4841309.233612698
$ws.Cells.Item(20, 8).Value2 = 3985931.394230286
$ws.Cells.Item(20, 9).Value2 = 1115523.345387527
$ws.Cells.Item(20, 10).Value2 = 4842344.303415114
$ws.Cells.Item(20, 11).Value2 = 3987525.566306008
$ws.Cells.Item(21, 6).Value2 = 1116430.340062525
$ws.Cells.Item(21, 7).Value2 = 4841325.746442241
$ws.Cells.Item(21, 8).Value2 = 3985957.360737001
$ws.Cells.Item(21, 9).Value2 = 1115569.229206563
$ws.Cells.Item(21, 10).Value2 = 4842295.651859458
$ws.Cells.Item(21, 11).Value2 = 3987556.686441578
$ws.Cells.Item(22, 6).Value2 = 1116427.564049069
$ws.Cells.Item(22, 7).Value2 = 4841342.259271783
$ws.Cells.Item(22, 8).Value2 = 3985981.25765702
$ws.Cells.Item(22, 9).Value2 = 1115616.242872429
$ws.Cells.Item(22, 10).Value2 = 4842247.000303801
$ws.Cells.Item(22, 11).Value2 = 3987572.657898724
$ws.Cells.Item(23, 6).Value2 = 1116425.089658014
$ws.Cells.Item(23, 7).Value2 = 4841358.772101326
$ws.Cells.Item(23, 8).Value2 = 3986003.390677121
$ws.Cells.Item(23, 9).Value2 = 1115664.414206565
$ws.Cells.Item(23, 10).Value2 = 4842198.348748145
$ws.Cells.Item(23, 11).Value2 = 3987573.480677446
$ws.Cells.Item(24, 6).Value2 = 1116422.866026873
$ws.Cells.Item(24, 7).Value2 = 4841375.284930869
$ws.Cells.Item(24, 8).Value2 = 3986024.002390224
$ws.Cells.Item(24, 9).Value2 = 1115713.771715488
$ws.Cells.Item(24, 10).Value2 = 4842149.697192488
$ws.Cells.Item(24, 11).Value2 = 3987559.154777746
$ws.Cells.Item(25, 6).Value2 = 1116420.853528542
$ws.Cells.Item(25, 7).Value2 = 4841391.797760411
$ws.Cells.Item(25, 8).Value2 = 3986043.288546725
$ws.Cells.Item(25, 9).Value2 = 1115764.34460766
$ws.Cells.Item(25, 10).Value2 = 4842101.045636832
$ws.Cells.Item(25, 11).Value2 = 3987529.680199622
$ws.Cells.Item(26, 6).Value2 = 1116419.020773035
$ws.Cells.Item(26, 7).Value2 = 4841408.310589955
$ws.Cells.Item(26, 8).Value2 = 3986061.409388488
$ws.Cells.Item(26, 9).Value2 = 1115816.162810775
$ws.Cells.Item(26, 10).Value2 = 4842052.394081175
$ws.Cells.Item(26, 11).Value2 = 3987485.056943076
$ws.Cells.Item(27, 6).Value2 = 1116417.342531475
$ws.Cells.Item(27, 7).Value2 = 4841424.823419497
$ws.Cells.Item(27, 8).Value2 = 3986078.497749532
$ws.Cells.Item(27, 9).Value2 = 1115869.25698947
$ws.Cells.Item(27, 10).Value2 = 4842003.742525518
$ws.Cells.Item(27, 11).Value2 = 3987425.285008106
$ws.Cells.Item(28, 6).Value2 = 1116415.798265123
$ws.Cells.Item(28, 7).Value2 = 4841441.33624904
$ws.Cells.Item(28, 8).Value2 = 3986094.664970791
$ws.Cells.Item(28, 9).Value2 = 1115923.65856347
$ws.Cells.Item(28, 10).Value2 = 4841955.090969862
$ws.Cells.Item(28, 11).Value2 = 3987350.364394713
$ws.Cells.Item(29, 6).Value2 = 1116414.371061704
$ws.Cells.Item(29, 7).Value2 = 4841457.849078584
$ws.Cells.Item(29, 8).Value2 = 3986110.005300671
$ws.Cells.Item(29, 9).Value2 = 1115979.399726181
$ws.Cells.Item(29, 10).Value2 = 4841906.439414206
$ws.Cells.Item(29, 11).Value2 = 3987260.295102897
$ws.Cells.Item(30, 6).Value2 = 1116413.046852223
$ws.Cells.Item(30, 7).Value2 = 4841474.361908126
$ws.Cells.Item(30, 8).Value2 = 3986124.599224009
$ws.Cells.Item(30, 9).Value2 = 1116036.513463743
$ws.Cells.Item(30, 10).Value2 = 4841857.78785855
$ws.Cells.Item(30, 11).Value2 = 3987155.077132658
$ws.Cells.Item(31, 6).Value2 = 1116411.813824913
$ws.Cells.Item(31, 7).Value2 = 4841490.874737669
$ws.Cells.Item(31, 8).Value2 = 3986138.516018078
$ws.Cells.Item(31, 9).Value2 = 1116095.033574549
$ws.Cells.Item(31, 10).Value2 = 4841809.136302892
$ws.Cells.Item(31, 11).Value2 = 3987034.710483995
$ws.Cells.Item(32, 6).Value2 = 1116410.66198028
$ws.Cells.Item(32, 7).Value2 = 4841507.387567212
$ws.Cells.Item(32, 8).Value2 = 3986151.815741494
$ws.Cells.Item(32, 9).Value2 = 1116154.994689246
$ws.Cells.Item(32, 10).Value2 = 4841760.484747237
$ws.Cells.Item(32, 11).Value2 = 3986899.195156909
$ws.Cells.Item(33, 6).Value2 = 1116409.58278879
$ws.Cells.Item(33, 7).Value2 = 4841523.900396754
$ws.Cells.Item(33, 8).Value2 = 3986164.550800603
$ws.Cells.Item(33, 9).Value2 = 1116216.432291229
$ws.Cells.Item(33, 10).Value2 = 4841711.83319158
$ws.Cells.Item(33, 11).Value2 = 3986748.5311514
$ws.Cells.Item(34, 6).Value2 = 1116408.568924368
$ws.Cells.Item(34, 7).Value2 = 4841540.413226298
$ws.Cells.Item(34, 8).Value2 = 3986176.767196691
$ws.Cells.Item(34, 9).Value2 = 1116279.382737639
$ws.Cells.Item(34, 10).Value2 = 4841663.181635924
$ws.Cells.Item(34, 11).Value2 = 3986582.718467468
$ws.Cells.Item(35, 6).Value2 = 1116407.614054636
$ws.Cells.Item(35, 7).Value2 = 4841556.92605584
$ws.Cells.Item(35, 8).Value2 = 3986188.505529017
$ws.Cells.Item(35, 9).Value2 = 1116343.883280879
$ws.Cells.Item(35, 10).Value2 = 4841614.530080266
$ws.Cells.Item(35, 11).Value2 = 3986401.757105112
$ws.Cells.Item(36, 6).Value2 = 1116406.712674155
$ws.Cells.Item(36, 7).Value2 = 4841573.438885383
$ws.Cells.Item(36, 8).Value2 = 3986199.801808876
$ws.Cells.Item(36, 9).Value2 = 1116409.972090656
$ws.Cells.Item(36, 10).Value2 = 4841565.87852461
$ws.Cells.Item(36, 11).Value2 = 3986205.647064334
$ws.Cells.Item(37, 6).Value2 = 1116405.859970621
$ws.Cells.Item(37, 7).Value2 = 4841589.951714926
$ws.Cells.Item(37, 8).Value2 = 3986210.688125865
$ws.Cells.Item(37, 9).Value2 = 1116477.688276574
$ws.Cells.Item(37, 10).Value2 = 4841517.226968954
$ws.Cells.Item(37, 11).Value2 = 3985994.388345132
$ws.Cells.Item(38, 6).Value2 = 1116405.051716608
$ws.Cells.Item(38, 7).Value2 = 4841606.464544469
$ws.Cells.Item(38, 8).Value2 = 3986221.193197436
$ws.Cells.Item(38, 9).Value2 = 1116547.071911273
$ws.Cells.Item(38, 10).Value2 = 4841468.575413298
$ws.Cells.Item(38, 11).Value2 = 3985767.980947508
$ws.Cells.Item(39, 6).Value2 = 1116404.284181261
$ws.Cells.Item(39, 7).Value2 = 4841622.977374012
$ws.Cells.Item(39, 8).Value2 = 3986231.34282543
$ws.Cells.Item(39, 9).Value2 = 1116618.164054144
$ws.Cells.Item(39, 10).Value2 = 4841419.92385764
$ws.Cells.Item(39, 11).Value2 = 3985526.424871459
$ws.Cells.Item(40, 6).Value2 = 1116403.554057765
$ws.Cells.Item(40, 7).Value2 = 4841639.490203555
$ws.Cells.Item(40, 8).Value2 = 3986241.160277881
$ws.Cells.Item(40, 9).Value2 = 1116691.006775629
$ws.Cells.Item(40, 10).Value2 = 4841371.272301984
$ws.Cells.Item(40, 11).Value2 = 3985269.720116988
$ws.Cells.Item(41, 6).Value2 = 1116402.858403341
$ws.Cells.Item(41, 7).Value2 = 4841656.003033097
$ws.Cells.Item(41, 8).Value2 = 3986250.666610292
$ws.Cells.Item(41, 9).Value2 = 1116765.643182117
$ws.Cells.Item(41, 10).Value2 = 4841322.620746328
$ws.Cells.Item(41, 11).Value2 = 3984997.866684094
$ws.Cells.Item(42, 6).Value2 = 1116402.194589309
$ws.Cells.Item(42, 7).Value2 = 4841672.515862641
$ws.Cells.Item(42, 8).Value2 = 3986259.880937555
$ws.Cells.Item(42, 9).Value2 = 1116842.117441451
$ws.Cells.Item(42, 10).Value2 = 4841273.969190672
$ws.Cells.Item(42, 11).Value2 = 3984710.864572776
$ws.Cells.Item(43, 6).Value2 = 1116401.560259268
$ws.Cells.Item(43, 7).Value2 = 4841689.028692184
$ws.Cells.Item(43, 8).Value2 = 3986268.820665341
$ws.Cells.Item(43, 9).Value2 = 1116920.474809066
$ws.Cells.Item(43, 10).Value2 = 4841225.317635015
$ws.Cells.Item(43, 11).Value2 = 3984408.713783035
$ws.Cells.Item(44, 6).Value2 = 1116400.953293884
$ws.Cells.Item(44, 7).Value2 = 4841705.541521726
$ws.Cells.Item(44, 8).Value2 = 3986277.501687993
$ws.Cells.Item(44, 9).Value2 = 1117000.76165477
$ws.Cells.Item(44, 10).Value2 = 4841176.666079358
$ws.Cells.Item(44, 11).Value2 = 3984091.414314871
$ws.Cells.Item(45, 6).Value2 = 1116400.371781094
$ws.Cells.Item(45, 7).Value2 = 4841722.054351269
$ws.Cells.Item(45, 8).Value2 = 3986285.938558591
$ws.Cells.Item(45, 9).Value2 = 1117083.025490188
$ws.Cells.Item(45, 10).Value2 = 4841128.014523703
$ws.Cells.Item(45, 11).Value2 = 3983758.966168284
$ws.Cells.Item(46, 6).Value2 = 1116399.813990753
$ws.Cells.Item(46, 7).Value2 = 4841738.567180811
$ws.Cells.Item(46, 8).Value2 = 3986294.14463574
$ws.Cells.Item(46, 9).Value2 = 1117167.314996872
$ws.Cells.Item(46, 10).Value2 = 4841079.362968045
$ws.Cells.Item(46, 11).Value2 = 3983411.369343274
$ws.Cells.Item(47, 6).Value2 = 1116399.278352975
$ws.Cells.Item(47, 7).Value2 = 4841755.080010355
$ws.Cells.Item(47, 8).Value2 = 3986302.132210815
$ws.Cells.Item(47, 9).Value2 = 1117253.680055112
$ws.Cells.Item(47, 10).Value2 = 4841030.711412389
$ws.Cells.Item(47, 11).Value2 = 3983048.62383984
$ws.Cells.Item(48, 6).Value2 = 1116398.763439543
$ws.Cells.Item(48, 7).Value2 = 4841771.592839898
$ws.Cells.Item(48, 8).Value2 = 3986309.912618706
$ws.Cells.Item(48, 9).Value2 = 1117342.171773458
$ws.Cells.Item(48, 10).Value2 = 4840982.059856732
$ws.Cells.Item(48, 11).Value2 = 3982670.729657983
$ws.Cells.Item(49, 6).Value2 = 1116398.267947888
$ws.Cells.Item(49, 7).Value2 = 4841788.10566944
$ws.Cells.Item(49, 8).Value2 = 3986317.496334551
$ws.Cells.Item(49, 9).Value2 = 1117432.842518957
$ws.Cells.Item(49, 10).Value2 = 4840933.408301076
$ws.Cells.Item(49, 11).Value2 = 3982277.686797704
$ws.Cells.Item(50, 6).Value2 = 1116397.790687228
$ws.Cells.Item(50, 7).Value2 = 4841804.618498984
$ws.Cells.Item(50, 8).Value2 = 3986324.893058546
$ws.Cells.Item(50, 9).Value2 = 1117525.745948148
$ws.Cells.Item(50, 10).Value2 = 4840884.756745419
$ws.Cells.Item(50, 11).Value2 = 3981869.495259
$ws.Cells.Item(51, 6).Value2 = 1116397.330566531
$ws.Cells.Item(51, 7).Value2 = 4841821.131328527
$ws.Cells.Item(51, 8).Value2 = 3986332.111790551
$ws.Cells.Item(51, 9).Value2 = 1117620.937038813
$ws.Cells.Item(51, 10).Value2 = 4840836.105189763
$ws.Cells.Item(51, 11).Value2 = 3981446.155041874
$ws.Cells.Item(52, 6).Value2 = 1116396.886584023
$ws.Cells.Item(52, 7).Value2 = 4841837.644158069
$ws.Cells.Item(52, 8).Value2 = 3986339.160895925
$ws.Cells.Item(52, 9).Value2 = 1117718.472122512
$ws.Cells.Item(52, 10).Value2 = 4840787.453634107
$ws.Cells.Item(52, 11).Value2 = 3981007.666146324
$ws.Cells.Item(53, 6).Value2 = 1116396.457818017
$ws.Cells.Item(53, 7).Value2 = 4841854.156987612
$ws.Cells.Item(53, 8).Value2 = 3986346.048163801
$ws.Cells.Item(53, 9).Value2 = 1117818.408917916
$ws.Cells.Item(53, 10).Value2 = 4840738.80207845
$ws.Cells.Item(53, 11).Value2 = 3980554.028572352
$ws.Cells.Item(54, 6).Value2 = 1116396.043418862
$ws.Cells.Item(54, 7).Value2 = 4841870.669817155
$ws.Cells.Item(54, 8).Value2 = 3986352.780858826
$ws.Cells.Item(54, 9).Value2 = 1117920.806564967
$ws.Cells.Item(54, 10).Value2 = 4840690.150522794
$ws.Cells.Item(54, 11).Value2 = 3980085.242319956
$ws.Cells.Item(55, 6).Value2 = 1116579.2486524
$ws.Cells.Item(55, 7).Value2 = 4841125.57814165
$ws.Cells.Item(55, 8).Value2 = 3985226.5088916
$ws.Cells.Item(55, 9).Value2 = 1114860.532467779
$ws.Cells.Item(55, 10).Value2 = 4843228.355522587
$ws.Cells.Item(55, 11).Value2 = 3984365.689238996
$ws.Cells.Item(56, 6).Value2 = 1116579.2486524
$ws.Cells.Item(56, 7).Value2 = 4841125.57814165
$ws.Cells.Item(56, 8).Value2 = 3985226.5088916
$ws.Cells.Item(56, 9).Value2 = 1114890.146442329
$ws.Cells.Item(56, 10).Value2 = 4843179.703883313
$ws.Cells.Item(56, 11).Value2 = 3984669.484877801
$ws.Cells.Item(57, 6).Value2 = 1116579.2486524
$ws.Cells.Item(57, 7).Value2 = 4841125.57814165
$ws.Cells.Item(57, 8).Value2 = 3985226.5088916
$ws.Cells.Item(57, 9).Value2 = 1114920.489633759
$ws.Cells.Item(57, 10).Value2 = 4843131.052244037
$ws.Cells.Item(57, 11).Value2 = 3984958.131873507
$ws.Cells.Item(58, 6).Value2 = 1116579.2486524
$ws.Cells.Item(58, 7).Value2 = 4841125.57814165
$ws.Cells.Item(58, 8).Value2 = 3985226.5088916
$ws.Cells.Item(58, 9).Value2 = 1114951.579998362
$ws.Cells.Item(58, 10).Value2 = 4843082.400604763
$ws.Cells.Item(58, 11).Value2 = 3985231.630226112
$ws.Cells.Item(59, 6).Value2 = 1116579.2486524
$ws.Cells.Item(59, 7).Value2 = 4841125.57814165
$ws.Cells.Item(59, 8).Value2 = 3985226.5088916
$ws.Cells.Item(59, 9).Value2 = 1114983.435934591
$ws.Cells.Item(59, 10).Value2 = 4843033.748965489
$ws.Cells.Item(59, 11).Value2 = 3985489.979935618
$ws.Cells.Item(60, 6).Value2 = 1116579.2486524
$ws.Cells.Item(60, 7).Value2 = 4841125.57814165
$ws.Cells.Item(60, 8).Value2 = 3985226.5088916
$ws.Cells.Item(60, 9).Value2 = 1115016.076293943
$ws.Cells.Item(60, 10).Value2 = 4842985.097326214
$ws.Cells.Item(60, 11).Value2 = 3985733.181002024
$ws.Cells.Item(61, 6).Value2 = 1116579.2486524
$ws.Cells.Item(61, 7).Value2 = 4841125.57814165
$ws.Cells.Item(61, 8).Value2 = 3985226.5088916
$ws.Cells.Item(61, 9).Value2 = 1115049.520392115
$ws.Cells.Item(61, 10).Value2 = 4842936.445686939
$ws.Cells.Item(61, 11).Value2 = 3985961.233425329
$ws.Cells.Item(62, 6).Value2 = 1116579.2486524
$ws.Cells.Item(62, 7).Value2 = 4841125.57814165
$ws.Cells.Item(62, 8).Value2 = 3985226.5088916
$ws.Cells.Item(62, 9).Value2 = 1115083.788020436
$ws.Cells.Item(62, 10).Value2 = 4842887.794047665
$ws.Cells.Item(62, 11).Value2 = 3986174.137205534
$ws.Cells.Item(63, 6).Value2 = 1116579.2486524
$ws.Cells.Item(63, 7).Value2 = 4841142.090964322
$ws.Cells.Item(63, 8).Value2 = 3985226.5088916
$ws.Cells.Item(63, 9).Value2 = 1115118.899457579
$ws.Cells.Item(63, 10).Value2 = 4842839.142408391
$ws.Cells.Item(63, 11).Value2 = 3986371.892342639
$ws.Cells.Item(64, 6).Value2 = 1116530.568056158
$ws.Cells.Item(64, 7).Value2 = 4841158.603786994
$ws.Cells.Item(64, 8).Value2 = 3985425.763207404
$ws.Cells.Item(64, 9).Value2 = 1115154.875481561
$ws.Cells.Item(64, 10).Value2 = 4842790.490769116
$ws.Cells.Item(64, 11).Value2 = 3986554.498836644
$ws.Cells.Item(65, 6).Value2 = 1116501.970647948
$ws.Cells.Item(65, 7).Value2 = 4841175.116609666
$ws.Cells.Item(65, 8).Value2 = 3985544.44926815
$ws.Cells.Item(65, 9).Value2 = 1115191.737382035
$ws.Cells.Item(65, 10).Value2 = 4842741.839129841
$ws.Cells.Item(65, 11).Value2 = 3986721.956687549
$ws.Cells.Item(66, 6).Value2 = 1116483.917865874
$ws.Cells.Item(66, 7).Value2 = 4841191.629432338
$ws.Cells.Item(66, 8).Value2 = 3985629.266174081
$ws.Cells.Item(66, 9).Value2 = 1115229.506972896
$ws.Cells.Item(66, 10).Value2 = 4842693.187490567
$ws.Cells.Item(66, 11).Value2 = 3986874.265895354
$ws.Cells.Item(67, 6).Value2 = 1116471.327928937
$ws.Cells.Item(67, 7).Value2 = 4841208.14225501
$ws.Cells.Item(67, 8).Value2 = 3985695.313053288
$ws.Cells.Item(67, 9).Value2 = 1115268.206605185
$ws.Cells.Item(67, 10).Value2 = 4842644.535851292
$ws.Cells.Item(67, 11).Value2 = 3987011.426460059
$ws.Cells.Item(68, 6).Value2 = 1116461.935458955
$ws.Cells.Item(68, 7).Value2 = 4841224.655077682
$ws.Cells.Item(68, 8).Value2 = 3985749.410570423
$ws.Cells.Item(68, 9).Value2 = 1115307.859180316
$ws.Cells.Item(68, 10).Value2 = 4842595.884212018
$ws.Cells.Item(68, 11).Value2 = 3987133.438381664
$ws.Cells.Item(69, 6).Value2 = 1116454.59185905
$ws.Cells.Item(69, 7).Value2 = 4841241.167900354
$ws.Cells.Item(69, 8).Value2 = 3985795.227218507
$ws.Cells.Item(69, 9).Value2 = 1115348.48816363
$ws.Cells.Item(69, 10).Value2 = 4842547.232572743
$ws.Cells.Item(69, 11).Value2 = 3987240.301660169
$ws.Cells.Item(70, 6).Value2 = 1116448.650352021
$ws.Cells.Item(70, 7).Value2 = 4841257.680723026
$ws.Cells.Item(70, 8).Value2 = 3985834.964816266
$ws.Cells.Item(70, 9).Value2 = 1115390.117598279
$ws.Cells.Item(70, 10).Value2 = 4842498.580933468
$ws.Cells.Item(70, 11).Value2 = 3987332.016295573
$ws.Cells.Item(71, 6).Value2 = 1116443.716784109
$ws.Cells.Item(71, 7).Value2 = 4841274.193545698
$ws.Cells.Item(71, 8).Value2 = 3985870.049182904
$ws.Cells.Item(71, 9).Value2 = 1115432.772119459
$ws.Cells.Item(71, 10).Value2 = 4842449.929294194
$ws.Cells.Item(71, 11).Value2 = 3987408.582287878
$ws.Cells.Item(72, 6).Value2 = 1116439.535956509
$ws.Cells.Item(72, 7).Value2 = 4841290.70636837
$ws.Cells.Item(72, 8).Value2 = 3985901.456751034
$ws.Cells.Item(72, 9).Value2 = 1115476.476968981
$ws.Cells.Item(72, 10).Value2 = 4842401.27765492
$ws.Cells.Item(72, 11).Value2 = 3987469.999637083
$ws.Cells.Item(73, 6).Value2 = 1116435.934581415
$ws.Cells.Item(73, 7).Value2 = 4841307.219191043
$ws.Cells.Item(73, 8).Value2 = 3985929.88558842
$ws.Cells.Item(73, 9).Value2 = 1115521.258010215
$ws.Cells.Item(73, 10).Value2 = 4842352.626015645
$ws.Cells.Item(73, 11).Value2 = 3987516.268343187
$ws.Cells.Item(74, 6).Value2 = 1116432.790394868
$ws.Cells.Item(74, 7).Value2 = 4841323.732013714
$ws.Cells.Item(74, 8).Value2 = 3985955.852085307
$ws.Cells.Item(74, 9).Value2 = 1115567.141743392
$ws.Cells.Item(74, 10).Value2 = 4842303.97437637
$ws.Cells.Item(74, 11).Value2 = 3987547.388406192
$ws.Cells.Item(75, 6).Value2 = 1116430.014375319
$ws.Cells.Item(75, 7).Value2 = 4841340.244836386
$ws.Cells.Item(75, 8).Value2 = 3985979.748996281
$ws.Cells.Item(75, 9).Value2 = 1115614.155321286
$ws.Cells.Item(75, 10).Value2 = 4842255.322737096
$ws.Cells.Item(75, 11).Value2 = 3987563.359826096
$ws.Cells.Item(76, 6).Value2 = 1116427.539978833
$ws.Cells.Item(76, 7).Value2 = 4841356.757659058
$ws.Cells.Item(76, 8).Value2 = 3986001.882008004
$ws.Cells.Item(76, 9).Value2 = 1115662.326565284
$ws.Cells.Item(76, 10).Value2 = 4842206.671097822
$ws.Cells.Item(76, 11).Value2 = 3987564.182602901
$ws.Cells.Item(77, 6).Value2 = 1116425.316342812
$ws.Cells.Item(77, 7).Value2 = 4841373.270481731
$ws.Cells.Item(77, 8).Value2 = 3986022.493713306
$ws.Cells.Item(77, 9).Value2 = 1115711.683981848
$ws.Cells.Item(77, 10).Value2 = 4842158.019458546
$ws.Cells.Item(77, 11).Value2 = 3987549.856736605
$ws.Cells.Item(78, 6).Value2 = 1116423.303840064
$ws.Cells.Item(78, 7).Value2 = 4841389.783304403
$ws.Cells.Item(78, 8).Value2 = 3986041.779862508
$ws.Cells.Item(78, 9).Value2 = 1115762.256779388
$ws.Cells.Item(78, 10).Value2 = 4842109.367819272
$ws.Cells.Item(78, 11).Value2 = 3987520.382227209
$ws.Cells.Item(79, 6).Value2 = 1116421.471080535
$ws.Cells.Item(79, 7).Value2 = 4841406.296127074
$ws.Cells.Item(79, 8).Value2 = 3986059.900697412
$ws.Cells.Item(79, 9).Value2 = 1115814.07488554
$ws.Cells.Item(79, 10).Value2 = 4842060.716179998
$ws.Cells.Item(79, 11).Value2 = 3987475.759074713
$ws.Cells.Item(80, 6).Value2 = 1116419.792835291
$ws.Cells.Item(80, 7).Value2 = 4841422.808949746
$ws.Cells.Item(80, 8).Value2 = 3986076.989051988
$ws.Cells.Item(80, 9).Value2 = 1115867.168964885
$ws.Cells.Item(80, 10).Value2 = 4842012.064540722
$ws.Cells.Item(80, 11).Value2 = 3987415.987279117
$ws.Cells.Item(81, 6).Value2 = 1116418.24856555
$ws.Cells.Item(81, 7).Value2 = 4841439.321772419
$ws.Cells.Item(81, 8).Value2 = 3986093.156267128
$ws.Cells.Item(81, 9).Value2 = 1115921.570437089
$ws.Cells.Item(81, 10).Value2 = 4841963.412901448
$ws.Cells.Item(81, 11).Value2 = 3987341.066840421
$ws.Cells.Item(82, 6).Value2 = 1116416.821358998
$ws.Cells.Item(82, 7).Value2 = 4841455.834595091
$ws.Cells.Item(82, 8).Value2 = 3986108.496591202
$ws.Cells.Item(82, 9).Value2 = 1115977.311495496
$ws.Cells.Item(82, 10).Value2 = 4841914.761262174
$ws.Cells.Item(82, 11).Value2 = 3987250.997758626
$ws.Cells.Item(83, 6).Value2 = 1116415.497146611
$ws.Cells.Item(83, 7).Value2 = 4841472.347417763
$ws.Cells.Item(83, 8).Value2 = 3986123.090509016
$ws.Cells.Item(83, 9).Value2 = 1116034.425126187
$ws.Cells.Item(83, 10).Value2 = 4841866.109622899
$ws.Cells.Item(83, 11).Value2 = 3987145.78003373
$ws.Cells.Item(84, 6).Value2 = 1116414.264116595
$ws.Cells.Item(84, 7).Value2 = 4841488.860240434
$ws.Cells.Item(84, 8).Value2 = 3986137.007297818
$ws.Cells.Item(84, 9).Value2 = 1116092.945127489
$ws.Cells.Item(84, 10).Value2 = 4841817.457983624
$ws.Cells.Item(84, 11).Value2 = 3987025.413665733
$ws.Cells.Item(85, 6).Value2 = 1116413.112269434
$ws.Cells.Item(85, 7).Value2 = 4841505.373063107
$ws.Cells.Item(85, 8).Value2 = 3986150.3070162
$ws.Cells.Item(85, 9).Value2 = 1116152.906129986
$ws.Cells.Item(85, 10).Value2 = 4841768.80634435
$ws.Cells.Item(85, 11).Value2 = 3986889.898654637
$ws.Cells.Item(86, 6).Value2 = 1116412.033075575
$ws.Cells.Item(86, 7).Value2 = 4841521.885885779
$ws.Cells.Item(86, 8).Value2 = 3986163.042070488
$ws.Cells.Item(86, 9).Value2 = 1116214.343617006
$ws.Cells.Item(86, 10).Value2 = 4841720.154705076
$ws.Cells.Item(86, 11).Value2 = 3986739.235000441
$ws.Cells.Item(87, 6).Value2 = 1116411.019208928
$ws.Cells.Item(87, 7).Value2 = 4841538.398708451
$ws.Cells.Item(87, 8).Value2 = 3986175.258461953
$ws.Cells.Item(87, 9).Value2 = 1116277.293945623
$ws.Cells.Item(87, 10).Value2 = 4841671.503065802
$ws.Cells.Item(87, 11).Value2 = 3986573.422703144
$ws.Cells.Item(88, 6).Value2 = 1116410.0643371
$ws.Cells.Item(88, 7).Value2 = 4841554.911531122
$ws.Cells.Item(88, 8).Value2 = 3986186.996789836
$ws.Cells.Item(88, 9).Value2 = 1116341.794368169
$ws.Cells.Item(88, 10).Value2 = 4841622.851426527
$ws.Cells.Item(88, 11).Value2 = 3986392.461762748
$ws.Cells.Item(89, 6).Value2 = 1116409.16295464
$ws.Cells.Item(89, 7).Value2 = 4841571.424353795
$ws.Cells.Item(89, 8).Value2 = 3986198.293065419
$ws.Cells.Item(89, 9).Value2 = 1116407.88305428
$ws.Cells.Item(89, 10).Value2 = 4841574.199787253
$ws.Cells.Item(89, 11).Value2 = 3986196.352179252
$ws.Cells.Item(90, 6).Value2 = 1116408.310249236
$ws.Cells.Item(90, 7).Value2 = 4841587.937176467
$ws.Cells.Item(90, 8).Value2 = 3986209.179378288
$ws.Cells.Item(90, 9).Value2 = 1116475.599113487
$ws.Cells.Item(90, 10).Value2 = 4841525.548147978
$ws.Cells.Item(90, 11).Value2 = 3985985.093952655
$ws.Cells.Item(91, 6).Value2 = 1116407.501993448
$ws.Cells.Item(91, 7).Value2 = 4841604.449999139
$ws.Cells.Item(91, 8).Value2 = 3986219.684445883
$ws.Cells.Item(91, 9).Value2 = 1116544.982618354
$ws.Cells.Item(91, 10).Value2 = 4841476.896508704
$ws.Cells.Item(91, 11).Value2 = 3985758.687082958
$ws.Cells.Item(92, 6).Value2 = 1116406.734456417
$ws.Cells.Item(92, 7).Value2 = 4841620.962821811
$ws.Cells.Item(92, 8).Value2 = 3986229.834070036
$ws.Cells.Item(92, 9).Value2 = 1116616.074628197
$ws.Cells.Item(92, 10).Value2 = 4841428.244869429
$ws.Cells.Item(92, 11).Value2 = 3985517.131570161
$ws.Cells.Item(93, 6).Value2 = 1116406.004331318
$ws.Cells.Item(93, 7).Value2 = 4841637.475644483
$ws.Cells.Item(93, 8).Value2 = 3986239.651518771
$ws.Cells.Item(93, 9).Value2 = 1116688.917213379
$ws.Cells.Item(93, 10).Value2 = 4841379.593230154
$ws.Cells.Item(93, 11).Value2 = 3985260.427414265
$ws.Cells.Item(94, 6).Value2 = 1116405.308675368
$ws.Cells.Item(94, 7).Value2 = 4841653.988467155
$ws.Cells.Item(94, 8).Value2 = 3986249.157847583
$ws.Cells.Item(94, 9).Value2 = 1116763.553480206
$ws.Cells.Item(94, 10).Value2 = 4841330.94159088
$ws.Cells.Item(94, 11).Value2 = 3984988.574615268
$ws.Cells.Item(95, 6).Value2 = 1116404.644859879
$ws.Cells.Item(95, 7).Value2 = 4841670.501289827
$ws.Cells.Item(95, 8).Value2 = 3986258.372171359
$ws.Cells.Item(95, 9).Value2 = 1116840.027596441
$ws.Cells.Item(95, 10).Value2 = 4841282.289951606
$ws.Cells.Item(95, 11).Value2 = 3984701.573173171
$ws.Cells.Item(96, 6).Value2 = 1116404.010528445
$ws.Cells.Item(96, 7).Value2 = 4841687.0141125
$ws.Cells.Item(96, 8).Value2 = 3986267.311895762
$ws.Cells.Item(96, 9).Value2 = 1116918.384817432
$ws.Cells.Item(96, 10).Value2 = 4841233.63831233
$ws.Cells.Item(96, 11).Value2 = 3984399.423087975
$ws.Cells.Item(97, 6).Value2 = 1116403.403561729
$ws.Cells.Item(97, 7).Value2 = 4841703.526935171
$ws.Cells.Item(97, 8).Value2 = 3986275.992915128
$ws.Cells.Item(97, 9).Value2 = 1116998.671512904
$ws.Cells.Item(97, 10).Value2 = 4841184.986673056
$ws.Cells.Item(97, 11).Value2 = 3984082.124359678
$ws.Cells.Item(98, 6).Value2 = 1116402.822047663
$ws.Cells.Item(98, 7).Value2 = 4841720.039757843
$ws.Cells.Item(98, 8).Value2 = 3986284.429782533
$ws.Cells.Item(98, 9).Value2 = 1117080.935194388
$ws.Cells.Item(98, 10).Value2 = 4841136.335033782
$ws.Cells.Item(98, 11).Value2 = 3983749.67698828
$ws.Cells.Item(99, 6).Value2 = 1116402.264256098
$ws.Cells.Item(99, 7).Value2 = 4841736.552580515
$ws.Cells.Item(99, 8).Value2 = 3986292.635856575
$ws.Cells.Item(99, 9).Value2 = 1117165.224543349
$ws.Cells.Item(99, 10).Value2 = 4841087.683394507
$ws.Cells.Item(99, 11).Value2 = 3983402.080973783
$ws.Cells.Item(100, 6).Value2 = 1116401.728617144
$ws.Cells.Item(100, 7).Value2 = 4841753.065403188
$ws.Cells.Item(100, 8).Value2 = 3986300.623428628
$ws.Cells.Item(100, 9).Value2 = 1117251.589439982
$ws.Cells.Item(100, 10).Value2 = 4841039.031755232
$ws.Cells.Item(100, 11).Value2 = 3983039.336316186
$ws.Cells.Item(101, 6).Value2 = 1116401.213702582
$ws.Cells.Item(101, 7).Value2 = 4841769.578225859
$ws.Cells.Item(101, 8).Value2 = 3986308.403833574
$ws.Cells.Item(101, 9).Value2 = 1117340.080992742
$ws.Cells.Item(101, 10).Value2 = 4840990.380115958
$ws.Cells.Item(101, 11).Value2 = 3982661.443015489
$ws.Cells.Item(102, 6).Value2 = 1116400.718209839
$ws.Cells.Item(102, 7).Value2 = 4841786.091048531
$ws.Cells.Item(102, 8).Value2 = 3986315.987546548
$ws.Cells.Item(102, 9).Value2 = 1117430.751568576
$ws.Cells.Item(102, 10).Value2 = 4840941.728476684
$ws.Cells.Item(102, 11).Value2 = 3982268.401071692
$ws.Cells.Item(103, 6).Value2 = 1116400.240948132
$ws.Cells.Item(103, 7).Value2 = 4841802.603871204
$ws.Cells.Item(103, 8).Value2 = 3986323.384267744
$ws.Cells.Item(103, 9).Value2 = 1117523.654823926
$ws.Cells.Item(103, 10).Value2 = 4840893.076837408
$ws.Cells.Item(103, 11).Value2 = 3981860.210484794
$ws.Cells.Item(104, 6).Value2 = 1116399.780826425
$ws.Cells.Item(104, 7).Value2 = 4841819.116693876
$ws.Cells.Item(104, 8).Value2 = 3986330.602997017
$ws.Cells.Item(104, 9).Value2 = 1117618.845736469
$ws.Cells.Item(104, 10).Value2 = 4840844.425198134
$ws.Cells.Item(104, 11).Value2 = 3981436.871254797
$ws.Cells.Item(105, 6).Value2 = 1116399.336842943
$ws.Cells.Item(105, 7).Value2 = 4841835.629516548
$ws.Cells.Item(105, 8).Value2 = 3986337.652099723
$ws.Cells.Item(105, 9).Value2 = 1117716.380637659
$ws.Cells.Item(105, 10).Value2 = 4840795.77355886
$ws.Cells.Item(105, 11).Value2 = 3980998.383381699
$ws.Cells.Item(106, 6).Value2 = 1116398.908075995
$ws.Cells.Item(106, 7).Value2 = 4841852.142339219
$ws.Cells.Item(106, 8).Value2 = 3986344.539364991
$ws.Cells.Item(106, 9).Value2 = 1117816.317246061
$ws.Cells.Item(106, 10).Value2 = 4840747.121919585
$ws.Cells.Item(106, 11).Value2 = 3980544.746865502
$ws.Cells.Item(107, 6).Value2 = 1116398.493675931
$ws.Cells.Item(107, 7).Value2 = 4841868.655161892
$ws.Cells.Item(107, 8).Value2 = 3986351.272057468
$ws.Cells.Item(107, 9).Value2 = 1117918.714701504
$ws.Cells.Item(107, 10).Value2 = 4840698.47028031
$ws.Cells.Item(107, 11).Value2 = 3980075.961706204

Write-Host "edit complete"
